# server/db.xlsx : bump the test row's "Preco_km" and "Email" sample values
# (Preco_km "1" -> "12", Email "a@a" -> "a@email") without disturbing any
# other cell content or formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D2 ("Preco_km") holds the numeric-looking text "1". Assigning a plain
# string to .Value would make Excel auto-coerce "12" into a real number,
# so force Text formatting just long enough to type it in, then drop the
# formatting override again so the cell's style stays exactly as it was.
$d2 = $ws.Range("D2")
$d2.NumberFormat = "@"
$d2.Value = "12"
$d2.ClearFormats()

# E2 ("Email") is plain (non-numeric) text, so a normal assignment is safe.
$ws.Range("E2").Value = "a@email"
